$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B45").Value = "       2019/6/20 20:30-21.30"
$ws.Range("C45").Value = " 写度量报告"
$ws.Range("C44").Copy()
$ws.Range("C45").PasteSpecial(-4122)
